# Apply updated market data values to Ixion_Profits sheets (H:N columns)
# as refreshed by the scheduled data-update runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1034.8572
$ws.Cells.Item(19, 9).Value = 625
$ws.Cells.Item(19, 10).Value = 1198.8
$ws.Cells.Item(19, 11).Value = 625
$ws.Cells.Item(19, 12).Value = 1198.8
$ws.Cells.Item(19, 13).Value = -450
$ws.Cells.Item(19, 14).Value = -1548.8

$ws.Cells.Item(38, 8).Value = 1217.1765
$ws.Cells.Item(38, 9).Value = 112
$ws.Cells.Item(38, 10).Value = 3243.3333
$ws.Cells.Item(38, 11).Value = 336
$ws.Cells.Item(38, 12).Value = 9729.999899999999
$ws.Cells.Item(38, 13).Value = 36
$ws.Cells.Item(38, 14).Value = -10473.9999

$ws.Cells.Item(58, 8).Value = 4423.2
$ws.Cells.Item(58, 10).Value = 5479
$ws.Cells.Item(58, 12).Value = 16437
$ws.Cells.Item(58, 14).Value = -16737

$ws.Cells.Item(64, 8).Value = 17000
$ws.Cells.Item(64, 9).Value = 17000
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 17000
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = -16752
$ws.Cells.Item(64, 14).Value = $null   # remove N64

$ws.Cells.Item(67, 8).Value = 17000
$ws.Cells.Item(67, 9).Value = 17000
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 11).Value = 17000
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 13).Value = -16142
$ws.Cells.Item(67, 14).Value = $null   # remove N67

$ws.Cells.Item(70, 8).Value = 2803.4
$ws.Cells.Item(70, 9).Value = 3320.4
$ws.Cells.Item(70, 11).Value = 9961.200000000001
$ws.Cells.Item(70, 13).Value = -9691.200000000001

$ws.Cells.Item(73, 8).Value = 2803.4
$ws.Cells.Item(73, 9).Value = 3320.4
$ws.Cells.Item(73, 11).Value = 9961.200000000001
$ws.Cells.Item(73, 13).Value = -9025.200000000001

$ws.Cells.Item(94, 8).Value = 2246.25
$ws.Cells.Item(94, 9).Value = 2246.25
$ws.Cells.Item(94, 11).Value = 2246.25
$ws.Cells.Item(94, 13).Value = -1795.25

$ws.Cells.Item(108, 8).Value = 59667
$ws.Cells.Item(108, 10).Value = 59667
$ws.Cells.Item(108, 12).Value = 59667
$ws.Cells.Item(108, 14).Value = -67347

$ws.Cells.Item(112, 8).Value = 1240.5714
$ws.Cells.Item(112, 10).Value = 1314
$ws.Cells.Item(112, 12).Value = 3942
$ws.Cells.Item(112, 14).Value = -6158

$ws.Cells.Item(115, 8).Value = 733
$ws.Cells.Item(115, 9).Value = 733
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 2199
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = -632
$ws.Cells.Item(115, 14).Value = $null   # remove N115

$ws.Cells.Item(118, 8).Value = 1018
$ws.Cells.Item(118, 9).Value = 789.5714
$ws.Cells.Item(118, 10).Value = 1337.8
$ws.Cells.Item(118, 11).Value = 2368.7142
$ws.Cells.Item(118, 12).Value = 4013.4
$ws.Cells.Item(118, 13).Value = -711.7142000000003
$ws.Cells.Item(118, 14).Value = -7327.4

$ws.Cells.Item(137, 8).Value = 1304.0613
$ws.Cells.Item(137, 9).Value = 1043.6578
$ws.Cells.Item(137, 10).Value = 2203.6365
$ws.Cells.Item(137, 11).Value = 3130.9734
$ws.Cells.Item(137, 12).Value = 6610.9095
$ws.Cells.Item(137, 13).Value = -580.9733999999999
$ws.Cells.Item(137, 14).Value = -11710.9095

$ws.Cells.Item(138, 8).Value = 3426.087
$ws.Cells.Item(138, 9).Value = 1538.72
$ws.Cells.Item(138, 10).Value = 5672.952
$ws.Cells.Item(138, 11).Value = 4616.16
$ws.Cells.Item(138, 12).Value = 17018.856
$ws.Cells.Item(138, 13).Value = 523.8400000000001
$ws.Cells.Item(138, 14).Value = -27298.856

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5821.375
$ws.Cells.Item(32, 9).Value = 4115.0894
$ws.Cells.Item(32, 11).Value = 4115.0894
$ws.Cells.Item(32, 13).Value = -3828.0894

$ws.Cells.Item(88, 8).Value = 2126.75
$ws.Cells.Item(88, 9).Value = 1800
$ws.Cells.Item(88, 10).Value = 2453.5
$ws.Cells.Item(88, 11).Value = 1800
$ws.Cells.Item(88, 12).Value = 2453.5
$ws.Cells.Item(88, 13).Value = -1394
$ws.Cells.Item(88, 14).Value = -3265.5

$ws.Cells.Item(91, 8).Value = 2126.75
$ws.Cells.Item(91, 9).Value = 1800
$ws.Cells.Item(91, 10).Value = 2453.5
$ws.Cells.Item(91, 11).Value = 1800
$ws.Cells.Item(91, 12).Value = 2453.5
$ws.Cells.Item(91, 13).Value = -396
$ws.Cells.Item(91, 14).Value = -5261.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1781.5385
$ws.Cells.Item(86, 9).Value = 1768.5714
$ws.Cells.Item(86, 10).Value = 1796.6666
$ws.Cells.Item(86, 11).Value = 1768.5714
$ws.Cells.Item(86, 12).Value = 1796.6666
$ws.Cells.Item(86, 13).Value = -645.5714
$ws.Cells.Item(86, 14).Value = -4042.6666

$ws.Cells.Item(89, 8).Value = 1781.5385
$ws.Cells.Item(89, 9).Value = 1768.5714
$ws.Cells.Item(89, 10).Value = 1796.6666
$ws.Cells.Item(89, 11).Value = 8842.857
$ws.Cells.Item(89, 12).Value = 8983.333000000001
$ws.Cells.Item(89, 13).Value = -3226.857
$ws.Cells.Item(89, 14).Value = -20215.333

$ws.Cells.Item(99, 8).Value = 250002130
$ws.Cells.Item(99, 9).Value = 333334500
$ws.Cells.Item(99, 10).Value = 5000
$ws.Cells.Item(99, 11).Value = 333334500
$ws.Cells.Item(99, 12).Value = 5000
$ws.Cells.Item(99, 13).Value = -333333002
$ws.Cells.Item(99, 14).Value = -7996

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 10428485
$ws.Cells.Item(99, 9).Value = 13417.714
$ws.Cells.Item(99, 11).Value = 13417.714
$ws.Cells.Item(99, 13).Value = -11919.714

$ws.Cells.Item(126, 8).Value = 10428485
$ws.Cells.Item(126, 9).Value = 13417.714
$ws.Cells.Item(126, 11).Value = 40253.142
$ws.Cells.Item(126, 13).Value = -37783.142

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 98902.53999999999
$ws.Cells.Item(5, 10).Value = 273897
$ws.Cells.Item(5, 12).Value = 821691
$ws.Cells.Item(5, 14).Value = -821915

$ws.Cells.Item(117, 8).Value = 23820932
$ws.Cells.Item(117, 9).Value = 14661.286
$ws.Cells.Item(117, 10).Value = 47627204
$ws.Cells.Item(117, 11).Value = 43983.858
$ws.Cells.Item(117, 12).Value = 142881612
$ws.Cells.Item(117, 13).Value = -40541.858
$ws.Cells.Item(117, 14).Value = -142888496

$ws.Cells.Item(121, 8).Value = 939.0303
$ws.Cells.Item(121, 10).Value = 1033.8572
$ws.Cells.Item(121, 12).Value = 3101.5716
$ws.Cells.Item(121, 14).Value = -5721.571599999999

$ws.Cells.Item(129, 8).Value = 19608770
$ws.Cells.Item(129, 9).Value = 25641602
$ws.Cells.Item(129, 10).Value = 2064
$ws.Cells.Item(129, 11).Value = 76924806
$ws.Cells.Item(129, 12).Value = 6192
$ws.Cells.Item(129, 13).Value = -76919806
$ws.Cells.Item(129, 14).Value = -16192

$ws.Cells.Item(135, 8).Value = 98902.53999999999
$ws.Cells.Item(135, 10).Value = 273897
$ws.Cells.Item(135, 12).Value = 2465073
$ws.Cells.Item(135, 14).Value = -2470143

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1239.5264
$ws.Cells.Item(102, 9).Value = 950.7857
$ws.Cells.Item(102, 10).Value = 2048
$ws.Cells.Item(102, 11).Value = 950.7857
$ws.Cells.Item(102, 12).Value = 2048
$ws.Cells.Item(102, 13).Value = 671.2143
$ws.Cells.Item(102, 14).Value = -5292

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 489.85715
$ws.Cells.Item(16, 9).Value = 489.8
$ws.Cells.Item(16, 10).Value = 490
$ws.Cells.Item(16, 11).Value = 489.8
$ws.Cells.Item(16, 12).Value = 490
$ws.Cells.Item(16, 13).Value = -319.8
$ws.Cells.Item(16, 14).Value = -830

$ws.Cells.Item(100, 8).Value = 1635.2941
$ws.Cells.Item(100, 9).Value = 1562.7693
$ws.Cells.Item(100, 10).Value = 1871
$ws.Cells.Item(100, 11).Value = 1562.7693
$ws.Cells.Item(100, 12).Value = 1871
$ws.Cells.Item(100, 13).Value = -1021.7693
$ws.Cells.Item(100, 14).Value = -2953
